$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to hold lat/long columns (A: latitude, B: longitude) across
# 4 rows. It now holds a single "addresses" column with 2 data rows, so drop
# column B entirely and drop row 4 entirely.
$ws.Range("B1:B4").EntireColumn.Delete()
$ws.Range("A4").EntireRow.Delete()

# Row 2 (A2) previously held a numeric latitude value with its own cell
# style (distinct from the header/other rows). Align its formatting with the
# rest of column A before/while overwriting its value, by copying A1's
# format onto it.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Write the new header + address values.
$ws.Range("A1").Value = "addresses"
$ws.Range("A2").Value = "dilsukhnagar,hyderabad"
$ws.Range("A3").Value = "stonehousepet, nellore"
